$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Apply formatting first (reusing/creating styles the same way Excel would) ----

# Build a one-off template cell holding the new "black text on white fill" style,
# then propagate it to the new detail rows (2-9) via copy/paste-special so that
# only a single new style entry is produced.
$tmpl = $ws.Range("G1")
$tmpl.ClearFormats()
$tmpl.Interior.Color = 16777215
$tmpl.Font.Color = 0
$tmpl.Copy()
$ws.Range("A2:E9").PasteSpecial(-4122)  # xlPasteFormats
$tmpl.Clear()

# Row 10 (new TOTAUX row) reuses the existing header style (bold white on gray).
$ws.Range("A1:E1").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)  # xlPasteFormats

# ---- Fill in values ----

$labels = @(
    "Chaudière 1 Défaut pompe ",
    "Chaudière 1 Commande pompe ",
    "Chaudière 1 Fin de course V2V ",
    "Chaudière 1 Commande V2V ",
    "Chaudière 2 Défaut pompe ",
    "Chaudière 2 Commande pompe ",
    "Chaudière 2 Fin de course V2V ",
    "Chaudière 2 Commande V2V "
)

# Per-row counts for columns B (Télé-Mesure), C (Télé-Signalisation), D (Télé-Réglage), E (Télé-Commande)
$counts = @(
    @(0,1,0,0),
    @(0,0,1,0),
    @(0,1,0,0),
    @(0,0,0,1),
    @(0,1,0,0),
    @(0,0,1,0),
    @(0,1,0,0),
    @(0,0,0,1)
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $labels[$i]
    $ws.Cells.Item($r, 2).Value = $counts[$i][0]
    $ws.Cells.Item($r, 3).Value = $counts[$i][1]
    $ws.Cells.Item($r, 4).Value = $counts[$i][2]
    $ws.Cells.Item($r, 5).Value = $counts[$i][3]
}

# TOTAUX row, now at row 10
$ws.Cells.Item(10, 1).Value = " TOTAUX (8 points)"
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 4
$ws.Cells.Item(10, 4).Value = 2
$ws.Cells.Item(10, 5).Value = 2
